$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet used to show area figures for three census years (1989, 2002,
# 2014) together with an extra census caption row. The refreshed export
# keeps only the most recent (2014) figure and drops the caption row, so:
#  1) remove row 2 ("(according to the population census data)")
#  2) remove the now-unused 1989/2002 data columns (B:C), leaving the
#     2014 column where column D used to be (it becomes column B)
$ws.Rows(2).Delete()
$ws.Columns("B:C").Delete()

# The refreshed template uses a taller uniform row height across the
# whole used area (including a couple of trailing blank rows).
for ($r = 1; $r -le 7; $r++) {
  $ws.Rows($r).RowHeight = 20.1
}
